$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.910.86"
$ws.Range("E2").Value = "  +2.27%  "

$ws.Range("D3").Value = "1.666.83"
$ws.Range("E3").Value = "  +1.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9984"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3649"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.59%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3243"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07075"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.52%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9982"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.077"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.13%  "

$ws.Range("D15").Value = "1.667.50"
$ws.Range("E15").Value = "  +1.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.611"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001052"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.91%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06622"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9986"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.935"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.22%  "

$ws.Range("D24").Value = "24.886.37"
$ws.Range("E24").Value = "  +2.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.446"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.419"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.58%  "

$ws.Range("D29").Value = "1.849.04"
$ws.Range("E29").Value = "  +1.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.02%  "

$ws.Range("E31").Value = "  +3.01%  "

$ws.Range("E32").Value = "  +0.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.747"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08489"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.642"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.166"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02258"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06036"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.224"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.247"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9978"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5939"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.839"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5676"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.67"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.953"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.186"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.79%  "
